$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text strings that look numeric.
# Force text format so Excel does not coerce them into real numbers
# (which would strip trailing zeros / use scientific notation).
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "272.19"
$ws.Range("D3").Value = "22.85"
$ws.Range("D4").Value = "6.470"
$ws.Range("D5").Value = "0.06222"
$ws.Range("D6").Value = "3.655"
$ws.Range("D7").Value = "6.658"
$ws.Range("D8").Value = "1.386"
$ws.Range("D9").Value = "0.8319"
$ws.Range("D10").Value = "0.01379"
$ws.Range("D11").Value = "0.1601"
$ws.Range("D12").Value = "0.08293"
$ws.Range("D13").Value = "0.03445"
$ws.Range("D14").Value = "0.03184"
$ws.Range("D15").Value = "0.09357"
$ws.Range("D16").Value = "3.849"
$ws.Range("D17").Value = "0.001641"
$ws.Range("D18").Value = "0.04733"
$ws.Range("D19").Value = "0.006292"
$ws.Range("D20").Value = "0.005690"
$ws.Range("D21").Value = "0.001075"
$ws.Range("D22").Value = "0.0001499"
$ws.Range("D23").Value = "3.718"
$ws.Range("D24").Value = "2.391"
$ws.Range("D25").Value = "0.3346"
$ws.Range("D27").Value = "0.0002703"
$ws.Range("D40").Value = "0.04700"
$ws.Range("D41").Value = "0.007034"
$ws.Range("D42").Value = "0.003797"
$ws.Range("D43").Value = "0.1161"
$ws.Range("D44").Value = "0.01168"
$ws.Range("D45").Value = "0.00006270"
$ws.Range("D48").Value = "0.9198"
$ws.Range("D49").Value = "0.00001399"
$ws.Range("D50").Value = "0.002099"
$ws.Range("D51").Value = "0.01239"

# Coin / link / label text columns
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJIWorstin24h"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("B49").Value = "CryptobidCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Range("E49").Value = "48CryptobidCoinCBC"
$ws.Range("B50").Value = "BOLO"
$ws.Range("C50").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("E50").Value = "49BOLOBOLOBestin24h"
